$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 50) with the latest ranking snapshot.
$ws.Range("A50").Value = "2025/12/04 19:00"
$ws.Range("B50").Value = "-"
$ws.Range("C50").Value = "-"
$ws.Range("D50").Value = "-"
$ws.Range("E50").Value = "-"
$ws.Range("F50").Value = "-"
$ws.Range("G50").Value = "-"
